# Updated symbol list on Mon Dec 26 00:55:32 UTC 2022 with GitHub Actions
# Refresh the crypto price table: new Price values, a couple of Volume(1h)
# labels, and roll the Data/Hora (date/hour) columns forward to 26-12-2022 / 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data (date) column: every data row moves from 25-12-2022 to 26-12-2022.
# This is plain text in the sheet (dd-mm-yyyy with dashes is not auto-parsed
# as a date by Excel), so a straight range assignment keeps it text.
$ws.Range("F2:F51").Value = "26-12-2022"

# --- Hora (hour) column: every data row moves from 23 to 0.
# Leading apostrophe forces Excel to store this as text (matching the
# original inline-string cell) instead of auto-converting "0" to a number.
$ws.Range("G2:G51").Value = "'0"

# --- Price column updates (text cells; force-text with a leading apostrophe
# so numeric-looking values like "242.72" stay text instead of becoming numbers).
$ws.Range("D2").Value = "'242.72"
$ws.Range("D3").Value = "'23.05"
$ws.Range("D4").Value = "'5.393"
$ws.Range("D5").Value = "'0.05986"
$ws.Range("D6").Value = "'3.402"
$ws.Range("D7").Value = "'6.484"
$ws.Range("D8").Value = "'0.8135"
$ws.Range("D9").Value = "'0.9007"
$ws.Range("D10").Value = "'0.1405"
$ws.Range("D11").Value = "'0.07411"
$ws.Range("D12").Value = "'0.03318"
$ws.Range("D13").Value = "'0.03071"
$ws.Range("D14").Value = "'0.09340"
$ws.Range("D15").Value = "'3.847"
$ws.Range("D16").Value = "'0.001565"
$ws.Range("D17").Value = "'0.04679"
$ws.Range("D18").Value = "'0.0005941"
$ws.Range("D19").Value = "'0.006107"
$ws.Range("D20").Value = "'0.005013"
$ws.Range("D21").Value = "'0.0009851"
$ws.Range("D22").Value = "'0.00007902"
$ws.Range("D23").Value = "'0.0002901"
$ws.Range("D25").Value = "'2.145"
$ws.Range("D26").Value = "'0.3205"
$ws.Range("D27").Value = "'0.1323"
$ws.Range("D40").Value = "'0.03894"
$ws.Range("D41").Value = "'0.006231"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D43").Value = "'0.002801"
$ws.Range("D44").Value = "'0.006493"
$ws.Range("D45").Value = "'0.00005205"
$ws.Range("D47").Value = "'0.0005801"
$ws.Range("D48").Value = "'0.9102"
$ws.Range("D49").Value = "'0.002299"
$ws.Range("D50").Value = "'0.00002101"

# --- Volume(1h) label updates (already non-numeric text, no force-text needed).
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

